$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.712.33"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "3.029.18"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'593.70"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'152.45"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.022.33"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "'6.56"
$ws.Range("E10").Value = "  +11.31%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "'35.49"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "3.530.40"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "'7.06"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "62.713.42"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "3.027.90"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "'451.45"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "'14.25"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").Value = "'83.05"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "'11.10"
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "'7.44"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  +7.18%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'27.52"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").Value = "0.0₃0863"
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = "  +8.70%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "'50.34"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.127"
$ws.Range("E41").Value = "  +3.42%  "
$ws.Range("D42").Value = "'9.07"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("E43").Value = "  +10.76%  "
$ws.Range("D44").Value = "'42.00"
$ws.Range("E44").Value = "  +4.74%  "
$ws.Range("D45").Value = "'393.24"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "2.720.10"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "'131.70"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +4.24%  "
$ws.Range("D51").Value = "'24.31"
$ws.Range("E51").Value = "  +3.58%  "
